$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 123
$ws.Range("I4").Value = 123
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 123
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -9
$ws.Range("H38").Value = 41.8
$ws.Range("I38").Value = 41.8
$ws.Range("K38").Value = 125.4
$ws.Range("M38").Value = 246.6
$ws.Range("H40").Value = 6343.077
$ws.Range("I40").Value = 5082.143
$ws.Range("K40").Value = 5082.143
$ws.Range("M40").Value = -4907.143
$ws.Range("H68").Value = 50295
$ws.Range("J68").Value = 50295
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51793
$ws.Range("H71").Value = 50295
$ws.Range("J71").Value = 50295
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158373
$ws.Range("H87").Value = 84087.60000000001
$ws.Range("J87").Value = 84087.60000000001
$ws.Range("L87").Value = 84087.60000000001
$ws.Range("N87").Value = -86583.60000000001
$ws.Range("H90").Value = 84087.60000000001
$ws.Range("J90").Value = 84087.60000000001
$ws.Range("L90").Value = 252262.8
$ws.Range("N90").Value = -264742.8
$ws.Range("H127").Value = 2062.125
$ws.Range("I127").Value = 1928.2858
$ws.Range("K127").Value = 5784.857400000001
$ws.Range("M127").Value = -824.8574000000008
$ws.Range("H129").Value = 977.8823
$ws.Range("I129").Value = 529
$ws.Range("K129").Value = 1587
$ws.Range("M129").Value = 3413
$ws.Range("H135").Value = 754.6539
$ws.Range("I135").Value = 651.4
$ws.Range("J135").Value = 1098.8334
$ws.Range("K135").Value = 5862.599999999999
$ws.Range("L135").Value = 9889.500599999999
$ws.Range("M135").Value = -3327.599999999999
$ws.Range("N135").Value = -14959.5006
$ws.Range("H138").Value = 3214.9443
$ws.Range("I138").Value = 1143.7693
$ws.Range("J138").Value = 8600
$ws.Range("K138").Value = 3431.3079
$ws.Range("L138").Value = 25800
$ws.Range("M138").Value = 1708.6921
$ws.Range("N138").Value = -36080
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 905.4286
$ws.Range("I2").Value = 905.4286
$ws.Range("K2").Value = 905.4286
$ws.Range("M2").Value = -792.4286
$ws.Range("H32").Value = 2780107
$ws.Range("I32").Value = 494.46667
$ws.Range("K32").Value = 494.46667
$ws.Range("M32").Value = -207.46667
$ws.Range("H74").Value = 3127.9
$ws.Range("I74").Value = 2733.4119
$ws.Range("J74").Value = 5363.3335
$ws.Range("K74").Value = 2733.4119
$ws.Range("L74").Value = 5363.3335
$ws.Range("M74").Value = -1859.4119
$ws.Range("N74").Value = -7111.3335
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H77").Value = 3127.9
$ws.Range("I77").Value = 2733.4119
$ws.Range("J77").Value = 5363.3335
$ws.Range("K77").Value = 13667.0595
$ws.Range("L77").Value = 26816.6675
$ws.Range("M77").Value = -9299.059499999999
$ws.Range("N77").Value = -35552.6675
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H80").Value = 51110
$ws.Range("J80").Value = 51110
$ws.Range("L80").Value = 51110
$ws.Range("N80").Value = -53106
$ws.Range("H83").Value = 51110
$ws.Range("J83").Value = 51110
$ws.Range("L83").Value = 153330
$ws.Range("N83").Value = -163314
$ws.Range("H116").Value = 905.4286
$ws.Range("I116").Value = 905.4286
$ws.Range("K116").Value = 905.4286
$ws.Range("M116").Value = 1388.5714
$ws.Range("H132").Value = 3728.2222
$ws.Range("I132").Value = 3653.4119
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10960.2357
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8430.235700000001
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 39494
$ws.Range("J134").Value = 39494
$ws.Range("L134").Value = 39494
$ws.Range("N134").Value = -49634
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 905.4286
$ws.Range("I3").Value = 905.4286
$ws.Range("K3").Value = 905.4286
$ws.Range("M3").Value = -791.4286
$ws.Range("H59").Value = 125780
$ws.Range("J59").Value = 125780
$ws.Range("L59").Value = 125780
$ws.Range("N59").Value = -127474
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H23").Value = 9000
$ws.Range("I23").Value = 8500
$ws.Range("K23").Value = 8500
$ws.Range("M23").Value = -8260
$ws.Range("H27").Value = 9000
$ws.Range("I27").Value = 8500
$ws.Range("K27").Value = 8500
$ws.Range("M27").Value = -8308
$ws.Range("H31").Value = 5336.8716
$ws.Range("I31").Value = 2343.8
$ws.Range("J31").Value = 6368.9653
$ws.Range("K31").Value = 2343.8
$ws.Range("L31").Value = 6368.9653
$ws.Range("M31").Value = -2048.8
$ws.Range("N31").Value = -6958.9653
$ws.Range("H32").Value = 1208.2727
$ws.Range("J32").Value = 2499.5
$ws.Range("L32").Value = 2499.5
$ws.Range("N32").Value = -3131.5
$ws.Range("H34").Value = 5336.8716
$ws.Range("I34").Value = 2343.8
$ws.Range("J34").Value = 6368.9653
$ws.Range("K34").Value = 2343.8
$ws.Range("L34").Value = 6368.9653
$ws.Range("M34").Value = -2141.8
$ws.Range("N34").Value = -6772.9653
$ws.Range("H58").Value = 2547.3125
$ws.Range("I58").Value = 1405.1538
$ws.Range("K58").Value = 1405.1538
$ws.Range("M58").Value = -1202.1538
$ws.Range("H105").Value = 1875
$ws.Range("I105").Value = 1750
$ws.Range("K105").Value = 1750
$ws.Range("M105").Value = -3
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0
$ws.Range("H136").Value = 2547.3125
$ws.Range("I136").Value = 1405.1538
$ws.Range("K136").Value = 4215.4614
$ws.Range("M136").Value = -1665.4614
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 26148.75
$ws.Range("I11").Value = 29384.285
$ws.Range("J11").Value = 3500
$ws.Range("K11").Value = 88152.855
$ws.Range("L11").Value = 10500
$ws.Range("M11").Value = -88012.855
$ws.Range("N11").Value = -10780
$ws.Range("H34").Value = 1712.2174
$ws.Range("J34").Value = 2352.5938
$ws.Range("L34").Value = 7057.7814
$ws.Range("N34").Value = -7225.7814
$ws.Range("H39").Value = 6759.8184
$ws.Range("J39").Value = 6759.8184
$ws.Range("L39").Value = 20279.4552
$ws.Range("N39").Value = -20867.4552
$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H55").Value = 3551.074
$ws.Range("J55").Value = 3668.8076
$ws.Range("L55").Value = 11006.4228
$ws.Range("N55").Value = -11360.4228
$ws.Range("H68").Value = 692
$ws.Range("I68").Value = 686.5
$ws.Range("K68").Value = 2059.5
$ws.Range("M68").Value = -1248.5
$ws.Range("H71").Value = 692
$ws.Range("I71").Value = 686.5
$ws.Range("K71").Value = 6178.5
$ws.Range("M71").Value = -2122.5
$ws.Range("H81").Value = 15
$ws.Range("J81").Value = 15
$ws.Range("L81").Value = 45
$ws.Range("N81").Value = -2291
$ws.Range("H84").Value = 15
$ws.Range("J84").Value = 15
$ws.Range("L84").Value = 135
$ws.Range("N84").Value = -11367
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.21429000000001
$ws.Range("I2").Value = 63.1
$ws.Range("K2").Value = 63.1
$ws.Range("M2").Value = 49.9
$ws.Range("H102").Value = 2383
$ws.Range("I102").Value = 2208.5334
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2208.5334
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -586.5333999999998
$ws.Range("N102").Value = -8244
$ws.Range("H126").Value = 3106.5715
$ws.Range("I126").Value = 3106.5715
$ws.Range("K126").Value = 9319.7145
$ws.Range("M126").Value = -6849.7145
$ws.Range("H132").Value = 28248.871
$ws.Range("I132").Value = 33594.25
$ws.Range("J132").Value = 3812.8572
$ws.Range("K132").Value = 100782.75
$ws.Range("L132").Value = 11438.5716
$ws.Range("M132").Value = -98252.75
$ws.Range("N132").Value = -16498.5716
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 999
$ws.Range("K16").Value = 999
$ws.Range("M16").Value = -829
$ws.Range("H61").Value = 66670400
$ws.Range("I61").Value = 125001120
$ws.Range("K61").Value = 125001120
$ws.Range("M61").Value = -125000918
$ws.Range("H113").Value = 66670400
$ws.Range("I113").Value = 125001120
$ws.Range("K113").Value = 125001120
$ws.Range("M113").Value = -124998950
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
